$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.961.13'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '3.154.81'
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = '216.79'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').Value = '625.53'
$ws.Range('E6').Value = '  +1.43%  '
$ws.Range('D7').Value = '1.13'
$ws.Range('E7').Value = '  +22.57%  '
$ws.Range('D8').Value = '0.368'
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = '3.153.65'
$ws.Range('E10').Value = '  +2.44%  '
$ws.Range('E11').Value = '  +11.91%  '
$ws.Range('E12').Value = '  +6.25%  '
$ws.Range('D13').Value = '5.71'
$ws.Range('E13').Value = '  +6.49%  '
$ws.Range('D14').Value = '0.0000246'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '35.10'
$ws.Range('E15').Value = '  +6.55%  '
$ws.Range('D16').Value = '90.640.13'
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '3.740.62'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D18').Value = '3.200.22'
$ws.Range('E18').Value = '  +4.00%  '
$ws.Range('D19').Value = '3.76'
$ws.Range('E19').Value = '  +8.56%  '
$ws.Range('D20').Value = '14.56'
$ws.Range('E20').Value = '  +5.74%  '
$ws.Range('D21').Value = '472.36'
$ws.Range('E21').Value = '  +8.73%  '
$ws.Range('D22').Value = '0.0000210'
$ws.Range('E22').Value = '  -4.25%  '
$ws.Range('D23').Value = '9.14'
$ws.Range('E23').Value = '  +8.10%  '
$ws.Range('D24').Value = '5.33'
$ws.Range('E24').Value = '  +3.71%  '
$ws.Range('D25').Value = '5.90'
$ws.Range('E25').Value = '  +4.63%  '
$ws.Range('D26').Value = '95.18'
$ws.Range('E26').Value = '  +13.30%  '
$ws.Range('E27').Value = '  +4.94%  '
$ws.Range('D28').Value = '3.321.50'
$ws.Range('E28').Value = '  +1.96%  '
$ws.Range('D30').Value = '0.236'
$ws.Range('E30').Value = '  +66.61%  '
$ws.Range('B31').Value = 'Cronos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D31').Value = '0.163'
$ws.Range('E31').Value = '  -2.26%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '9.31'
$ws.Range('E32').Value = '  +5.29%  '
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').Value = '27.14'
$ws.Range('E34').Value = '  +18.11%  '
$ws.Range('D35').Value = '519.72'
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('E36').Value = '  +5.48%  '
$ws.Range('E37').Value = '  +5.97%  '
$ws.Range('E38').Value = '  -0.80%  '
$ws.Range('E39').Value = '  -6.26%  '
$ws.Range('E40').Value = '  +3.36%  '
$ws.Range('D41').Value = '0.0923'
$ws.Range('E41').Value = '  +27.78%  '
$ws.Range('D42').Value = '0.430'
$ws.Range('E42').Value = '  +16.84%  '
$ws.Range('D43').Value = '22.24'
$ws.Range('E43').Value = '  -0.32%  '
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').Value = '1.98'
$ws.Range('E45').Value = '  +5.72%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '0.730'
$ws.Range('E47').Value = '  +20.52%  '
$ws.Range('E48').Value = '  +12.81%  '
$ws.Range('D49').Value = '150.07'
$ws.Range('E49').Value = '  +6.34%  '
$ws.Range('E50').Value = '  +10.07%  '
$ws.Range('D51').Value = '45.37'
$ws.Range('E51').Value = '  +3.46%  '
